$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title date line
Replace-Text "2025-03-15 Saturday" "2025-03-16 Sunday"

# Table cells with unique old values -> safe to use Find/Replace
Replace-Text "26÷2=13, 0" "34÷4=8, 2"
Replace-Text "54÷6=9, 0" "84÷6=14, 0"
Replace-Text "70÷9=7, 7" "73÷5=14, 3"
Replace-Text "37÷2=18, 1" "28÷9=3, 1"
Replace-Text "87÷7=12, 3" "17÷9=1, 8"

Replace-Text "93÷7=13, 2" "57÷2=28, 1"
Replace-Text "71÷3=23, 2" "45÷2=22, 1"
Replace-Text "28÷6=4, 4" "22÷7=3, 1"
Replace-Text "42÷2=21, 0" "72÷7=10, 2"
Replace-Text "35÷6=5, 5" "12÷4=3, 0"

Replace-Text "61÷3=20, 1" "60÷6=10, 0"
Replace-Text "27÷7=3, 6" "13÷9=1, 4"
Replace-Text "53÷2=26, 1" "43÷8=5, 3"
Replace-Text "98÷3=32, 2" "78÷9=8, 6"
Replace-Text "97÷8=12, 1" "11÷4=2, 3"

Replace-Text "27÷6=4, 3" "35÷9=3, 8"
Replace-Text "75÷7=10, 5" "57÷9=6, 3"
Replace-Text "68÷4=17, 0" "19÷6=3, 1"
# "32÷5=6, 2" appears twice in the same row (columns 4 and 5) and must
# become two different values, so address the cells directly instead of
# using Find/Replace (which would otherwise turn both into one value).
$t = $d.Tables.Item(1)
$t.Cell(13, 4).Range.Text = "23÷4=5, 3"
$t.Cell(13, 5).Range.Text = "58÷7=8, 2"

Replace-Text "24÷4=6, 0" "51÷5=10, 1"
Replace-Text "90÷9=10, 0" "60÷6=10, 0"
Replace-Text "57÷8=7, 1" "65÷3=21, 2"
Replace-Text "44÷8=5, 4" "63÷4=15, 3"
Replace-Text "62÷9=6, 8" "74÷6=12, 2"
